# Add two new rows (65, 66) of regression-test results to the
# "2010 and 2010-18" worksheet, matching the new
# "Demo_Baseline_2010-18_Dec22" / "Demo_Baseline_2010-18_Dec22_1800"
# scenario runs, and update the sheet view (frozen pane scroll spot and
# the active selection) to point at the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-NumCell($ws, $row, $col, $style, $value) {
    $cell = $ws.Cells.Item($row, $col)
    if ($style -eq 2) {
        $cell.NumberFormat = "0.00"
    } elseif ($style -eq 3) {
        $cell.NumberFormat = "0"
    } elseif ($style -eq 4) {
        $cell.NumberFormat = "0.000000"
    } elseif ($style -eq 5) {
        $cell.NumberFormat = "0.00"
        $cell.Interior.ColorIndex = 6
    } elseif ($style -eq 6) {
        $cell.NumberFormat = "0"
        $cell.Interior.ColorIndex = 6
    }
    $cell.Value2 = $value
}

# ---- Row 65: Demo_Baseline_2010-18_Dec22 ----
$ws.Cells.Item(65, 1).Value2 = "CW3M"
$ws.Cells.Item(65, 2).Value2 = "Demo_Baseline_2010-18_Dec22"
$ws.Cells.Item(65, 3).Value2 = "2010-18"

Set-NumCell $ws 65 4  5 1246.3303018888889
Set-NumCell $ws 65 5  5 1890.2624783333331
Set-NumCell $ws 65 6  2 0.94846033333333346
Set-NumCell $ws 65 7  2 305.6782124444444
Set-NumCell $ws 65 8  2 9.775355222222224
Set-NumCell $ws 65 9  5 6.3587768888888885
Set-NumCell $ws 65 10 2 8.145128999999999
Set-NumCell $ws 65 11 5 628.90788111111101
Set-NumCell $ws 65 12 5 82.308506444444433
Set-NumCell $ws 65 13 5 1456.3722873333334
Set-NumCell $ws 65 14 5 1283.7495253333334
Set-NumCell $ws 65 15 6 3986.0738390000001
Set-NumCell $ws 65 16 3 27227.338324888889
Set-NumCell $ws 65 17 2 0.12974411111111114
Set-NumCell $ws 65 18 4 0.000023444444444444448

$ws.Cells.Item(65, 19).Value2 = "2010-18"

# ---- Row 66: Demo_Baseline_2010-18_Dec22_1800 ----
$ws.Cells.Item(66, 1).Value2 = "CW3M"
$ws.Cells.Item(66, 2).Value2 = "Demo_Baseline_2010-18_Dec22_1800"
$ws.Cells.Item(66, 3).Value2 = "2010-18"

Set-NumCell $ws 66 4  5 1380.5085448888888
Set-NumCell $ws 66 5  2 1890.2624783333331
Set-NumCell $ws 66 6  2 0.94846033333333346
Set-NumCell $ws 66 7  2 305.6782124444444
Set-NumCell $ws 66 8  2 9.775355222222224
Set-NumCell $ws 66 9  2 6.1288343333333337
Set-NumCell $ws 66 10 2 8.145128999999999
Set-NumCell $ws 66 11 2 628.9703199999999
Set-NumCell $ws 66 12 2 82.308506444444433
Set-NumCell $ws 66 13 5 1487.9843207777781
Set-NumCell $ws 66 14 5 1386.0215385555557
Set-NumCell $ws 66 15 6 4013.3543294444444
Set-NumCell $ws 66 16 3 27227.338324888889
Set-NumCell $ws 66 17 2 0.12792911111111113
Set-NumCell $ws 66 18 4 0.000021888888888888887

# ---- Update frozen-pane scroll position & current selection ----
$win = $excel.ActiveWindow
$win.ScrollRow = 54
$win.ScrollColumn = 1
$ws.Range("M66:O66").Select()
